$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings in column D stay as literal text
# (matching the source data which stores prices/volumes as text strings,
# e.g. "76.538.36" or "0.550" with preserved trailing zeros).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "76.490.99"
$ws.Range("E2").Value = "  +0.54%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.940.75"
$ws.Range("E3").Value = "  +1.48%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "198.34"
$ws.Range("E5").Value = "  +0.77%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "594.46"
$ws.Range("E6").Value = "  -1.11%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.550"
$ws.Range("E8").Value = "  -1.13%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.200"
$ws.Range("E9").Value = "  +2.76%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "2.937.93"
$ws.Range("E10").Value = "  +1.43%  "
$ws.Range("E11").Value = "  +10.38%  "
$ws.Range("E12").Value = "  +0.35%  "
$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.480.69"
$ws.Range("E13").Value = "  +1.33%  "
$ws.Range("B14").Value = "Toncoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.89"
$ws.Range("E14").Value = "  -1.10%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "28.45"
$ws.Range("E15").Value = "  +2.85%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "76.501.20"
$ws.Range("E16").Value = "  +0.67%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000190"
$ws.Range("E17").Value = "  -0.95%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.929.08"
$ws.Range("E18").Value = "  +1.12%  "
$ws.Range("E19").Value = "  +6.95%  "
$ws.Range("E20").Value = "  -3.21%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "374.76"
$ws.Range("E21").Value = "  -2.40%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.33"
$ws.Range("E22").Value = "  +3.82%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.26"
$ws.Range("E23").Value = "  -2.58%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "71.90"
$ws.Range("E24").Value = "  -0.35%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  +0.09%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.092.97"
$ws.Range("E26").Value = "  +1.52%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "4.26"
$ws.Range("E27").Value = "  -0.48%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.66"
$ws.Range("E28").Value = "  -2.20%  "
$ws.Range("E29").Value = "  -1.09%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  +0.15%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.34"
$ws.Range("E31").Value = "  +6.08%  "
$ws.Range("E32").Value = "  -3.20%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "498.14"
$ws.Range("E33").Value = "  -3.04%  "
$ws.Range("E34").Value = "  -0.71%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "165.35"
$ws.Range("E36").Value = "  +0.18%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "20.15"
$ws.Range("E37").Value = "  -0.62%  "
$ws.Range("E38").Value = "  +12.39%  "
$ws.Range("E39").Value = "  +17.34%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "19.96"
$ws.Range("E40").Value = "  +1.33%  "
$ws.Range("E41").Value = "  -4.23%  "
$ws.Range("E42").Value = "  +0.03%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "179.59"
$ws.Range("E43").Value = "  -2.65%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.91"
$ws.Range("E44").Value = "  -3.36%  "
$ws.Range("E45").Value = "  -2.36%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "39.95"
$ws.Range("E46").Value = "  -0.91%  "
$ws.Range("E47").Value = "  -4.28%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.595"
$ws.Range("E48").Value = "  +0.83%  "
$ws.Range("E49").Value = "  -2.89%  "
$ws.Range("E50").Value = "  +2.09%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.665"
$ws.Range("E51").Value = "  -1.36%  "
